$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap developer names assigned to "LOAD DATA" (row 2) and "PREPROCESSING DATA" (row 3)
$e2 = $ws.Range("E2").Value()
$e3 = $ws.Range("E3").Value()
$ws.Range("E2").Value = $e3
$ws.Range("E3").Value = $e2

# Update the active selection to reflect the last cell clicked after editing
$ws.Activate()
$ws.Range("F2").Select()
